# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" data table (rows 16-59) was re-sorted
# from descending period order (2003 -> 1608) to ascending period order
# (1608 -> 2003). Read the current block, reverse the row order, and write
# it back so each period keeps the value it had before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 59
$periodCol = 5   # E - Periodo Mora
$valueCol = 6    # F - Valor Mora

$rowCount = $lastRow - $firstRow + 1

# Capture current values for the two columns that actually vary per row.
$periods = @()
$values = @()
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    $periods += $ws.Cells.Item($r, $periodCol).Value()
    $values += $ws.Cells.Item($r, $valueCol).Value()
}

# Write them back in reverse order.
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    $srcIndex = $rowCount - 1 - $i
    $ws.Cells.Item($r, $periodCol).Value = $periods[$srcIndex]
    $ws.Cells.Item($r, $valueCol).Value = $values[$srcIndex]
}
